$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = -8.435999999999998
$ws.Range("D6").Value = -7.837000000000002
$ws.Range("C7").Value = -13.463
$ws.Range("A8").Value = -21.255
$ws.Range("D9").Value = -7.939
$ws.Range("A10").Value = -20.945
$ws.Range("D10").Value = -7.673
$ws.Range("A12").Value = -21.808
$ws.Range("B13").Value = 6.475
$ws.Range("A18").Value = -21.78
$ws.Range("C20").Value = -13.041
